$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2987.8235
$ws.Range("I76").Value = 2580.6365
$ws.Range("J76").Value = 3099.8
$ws.Range("K76").Value = 2580.6365
$ws.Range("L76").Value = 3099.8
$ws.Range("M76").Value = -2265.6365
$ws.Range("N76").Value = -3729.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2987.8235
$ws.Range("I79").Value = 2580.6365
$ws.Range("J79").Value = 3099.8
$ws.Range("K79").Value = 2580.6365
$ws.Range("L79").Value = 3099.8
$ws.Range("M79").Value = -1488.6365
$ws.Range("N79").Value = -5283.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 51999.35
$ws.Range("I137").Value = 2125.9092
$ws.Range("K137").Value = 6377.7276
$ws.Range("M137").Value = -3827.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23813074
$ws.Range("I32").Value = 25001978
$ws.Range("J32").Value = 35000
$ws.Range("K32").Value = 25001978
$ws.Range("L32").Value = 35000
$ws.Range("M32").Value = -25001691
$ws.Range("N32").Value = -35574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3537.3333
$ws.Range("I61").Value = 3504.2856
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3504.2856
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3292.2856
$ws.Range("N61").Value = -4424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 507.52
$ws.Range("I97").Value = 481.64706
$ws.Range("J97").Value = 562.5
$ws.Range("K97").Value = 481.64706
$ws.Range("L97").Value = 562.5
$ws.Range("M97").Value = 14.35293999999999
$ws.Range("N97").Value = -1554.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3537.3333
$ws.Range("I136").Value = 3504.2856
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 10512.8568
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -7962.856800000001
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2472.158
$ws.Range("I99").Value = 1790
$ws.Range("J99").Value = 4382.2
$ws.Range("K99").Value = 1790
$ws.Range("L99").Value = 4382.2
$ws.Range("M99").Value = -292
$ws.Range("N99").Value = -7378.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 12657
$ws.Range("J103").Value = 12657
$ws.Range("L103").Value = 12657
$ws.Range("N103").Value = -15001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1258.0728
$ws.Range("I134").Value = 1272.26
$ws.Range("J134").Value = 1116.2
$ws.Range("K134").Value = 3816.78
$ws.Range("L134").Value = 3348.6
$ws.Range("M134").Value = -1281.78
$ws.Range("N134").Value = -8418.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23298.94
$ws.Range("I31").Value = 6743.769
$ws.Range("J31").Value = 34059.8
$ws.Range("K31").Value = 6743.769
$ws.Range("L31").Value = 34059.8
$ws.Range("M31").Value = -6448.769
$ws.Range("N31").Value = -34649.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 23298.94
$ws.Range("I34").Value = 6743.769
$ws.Range("J34").Value = 34059.8
$ws.Range("K34").Value = 6743.769
$ws.Range("L34").Value = 34059.8
$ws.Range("M34").Value = -6541.769
$ws.Range("N34").Value = -34463.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2386.5
$ws.Range("I58").Value = 1866.25
$ws.Range("J58").Value = 2733.3333
$ws.Range("K58").Value = 1866.25
$ws.Range("L58").Value = 2733.3333
$ws.Range("M58").Value = -1663.25
$ws.Range("N58").Value = -3139.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2626.3
$ws.Range("I99").Value = 2415.2144
$ws.Range("J99").Value = 3118.8333
$ws.Range("K99").Value = 2415.2144
$ws.Range("L99").Value = 3118.8333
$ws.Range("M99").Value = -917.2143999999998
$ws.Range("N99").Value = -6114.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2626.3
$ws.Range("I126").Value = 2415.2144
$ws.Range("J126").Value = 3118.8333
$ws.Range("K126").Value = 7245.6432
$ws.Range("L126").Value = 9356.499899999999
$ws.Range("M126").Value = -4775.6432
$ws.Range("N126").Value = -14296.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1898.2941
$ws.Range("I132").Value = 1448
$ws.Range("K132").Value = 4344
$ws.Range("M132").Value = -1814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3758.5
$ws.Range("I134").Value = 2199.4375
$ws.Range("J134").Value = 5837.25
$ws.Range("K134").Value = 6598.3125
$ws.Range("L134").Value = 17511.75
$ws.Range("M134").Value = -4063.3125
$ws.Range("N134").Value = -22581.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2386.5
$ws.Range("I136").Value = 1866.25
$ws.Range("J136").Value = 2733.3333
$ws.Range("K136").Value = 5598.75
$ws.Range("L136").Value = 8199.999899999999
$ws.Range("M136").Value = -3048.75
$ws.Range("N136").Value = -13299.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 447.9091
$ws.Range("I107").Value = 472.7143
$ws.Range("J107").Value = 404.5
$ws.Range("K107").Value = 1418.1429
$ws.Range("L107").Value = 1213.5
$ws.Range("M107").Value = 501.8571000000002
$ws.Range("N107").Value = -5053.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2842.2856
$ws.Range("I132").Value = 2321
$ws.Range("J132").Value = 5970
$ws.Range("K132").Value = 6963
$ws.Range("L132").Value = 17910
$ws.Range("M132").Value = -4433
$ws.Range("N132").Value = -22970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1730.7222
$ws.Range("I122").Value = 1746.0834
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 5238.2502
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -2788.2502
$ws.Range("N122").Value = -10000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3938.1667
$ws.Range("I132").Value = 4188.9
$ws.Range("J132").Value = 3624.75
$ws.Range("K132").Value = 12566.7
$ws.Range("L132").Value = 10874.25
$ws.Range("M132").Value = -10036.7
$ws.Range("N132").Value = -15934.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1975.359
$ws.Range("I122").Value = 1329.0476
$ws.Range("J122").Value = 2729.389
$ws.Range("K122").Value = 3987.142800000001
$ws.Range("L122").Value = 8188.167
$ws.Range("M122").Value = -1537.142800000001
$ws.Range("N122").Value = -13088.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3733.2
$ws.Range("I132").Value = 2999.9
$ws.Range("J132").Value = 5199.8
$ws.Range("K132").Value = 8999.700000000001
$ws.Range("L132").Value = 15599.4
$ws.Range("M132").Value = -6469.700000000001
$ws.Range("N132").Value = -20659.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20354.371
$ws.Range("I136").Value = 26353.45
$ws.Range("J136").Value = 3214.1428
$ws.Range("K136").Value = 79060.35000000001
$ws.Range("L136").Value = 9642.428400000001
$ws.Range("M136").Value = -14742.4284
